# Crypto list refresh (prices + 1h volume deltas), GitHub Actions run
# Wed Oct 18 20:01:56 UTC 2023. Row 51 coin swapped BabyDogeCoin -> Cronos.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.306.33"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.566.05"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'210.99"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'0.490"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'44.23"
$ws.Range("E8").Value = "  -4.21%  "
$ws.Range("D9").Value = "'23.75"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "'0.0894"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "1.792.85"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "1.568.08"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "'3.66"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "28.342.68"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "'0.512"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "'60.92"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").Value = "'227.16"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "0.0₃0680"
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("D23").Value = "'3.93"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "'8.95"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "'150.68"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "'14.90"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "'6.32"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +3.15%  "
$ws.Range("D32").Value = "'1.08"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").Value = "'3.07"
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("D35").Value = "1.379.14"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  +2.45%  "
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").Value = "'2.66"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'0.780"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D47").Value = "'62.25"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "1.705.11"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").Value = "'0.915"
$ws.Range("E49").Value = "  -6.39%  "
$ws.Range("D50").Value = "'85.38"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0514"
$ws.Range("E51").Value = "  -1.06%  "

# Re-normalize style on the text-forced Price cells so the apostrophe
# (quote-prefix) entry trick does not leave a lingering cell style change.
$textForcedRefs = @("D5","D6","D8","D9","D12","D15","D17","D18","D19","D23","D24","D26","D27","D29","D32","D34","D39","D45","D47","D49","D50","D51")
foreach ($ref in $textForcedRefs) {
    $ws.Range($ref).Style = "Normal"
}

Write-Host "Applied crypto list updates: 78 cells across 50 rows"
